$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '45.203.79'
Set-TextValue 'E2' '  +4.92%  '
Set-TextValue 'D3' '2.357.07'
Set-TextValue 'E3' '  +1.95%  '
Set-TextValue 'E4' '  +0.31%  '
Set-TextValue 'D5' '310.13'
Set-TextValue 'E5' '  -0.34%  '
Set-TextValue 'D6' '107.70'
Set-TextValue 'E6' '  +0.77%  '
Set-TextValue 'E7' '  -0.13%  '
Set-TextValue 'E8' '  -0.19%  '
Set-TextValue 'E9' '  +1.20%  '
Set-TextValue 'D10' '41.06'
Set-TextValue 'E10' '  +2.54%  '
Set-TextValue 'E11' '  +0.17%  '
Set-TextValue 'E12' '  +0.68%  '
Set-TextValue 'E13' '  +2.00%  '
Set-TextValue 'D14' '0.981'
Set-TextValue 'E14' '  -0.71%  '
Set-TextValue 'D15' '2.716.57'
Set-TextValue 'E15' '  +1.96%  '
Set-TextValue 'E16' '  +0.30%  '
Set-TextValue 'D17' '2.353.33'
Set-TextValue 'E17' '  +2.09%  '
Set-TextValue 'D18' '45.158.22'
Set-TextValue 'E18' '  +5.64%  '
Set-TextValue 'B19' 'InternetComputer(DFINITY)'
Set-TextValue 'C19' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D19' '14.02'
Set-TextValue 'E19' '  +7.24%  '
Set-TextValue 'B20' 'Uniswap'
Set-TextValue 'C20' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D20' '7.32'
Set-TextValue 'E20' '  -1.85%  '
Set-TextValue 'E21' '  +0.65%  '
Set-TextValue 'D22' '73.14'
Set-TextValue 'E22' '  -0.66%  '
Set-TextValue 'E23' '  -0.12%  '
Set-TextValue 'D24' '259.53'
Set-TextValue 'E24' '  -2.54%  '
Set-TextValue 'E25' '  +3.29%  '
Set-TextValue 'E26' '  -0.33%  '
Set-TextValue 'D27' '11.13'
Set-TextValue 'D28' '7.29'
Set-TextValue 'E28' '  -4.74%  '
Set-TextValue 'E29' '  +2.52%  '
Set-TextValue 'D30' '0.0968'
Set-TextValue 'E30' '  +10.53%  '
Set-TextValue 'D31' '22.29'
Set-TextValue 'E31' '  -0.75%  '
Set-TextValue 'D32' '37.74'
Set-TextValue 'E32' '  -2.56%  '
Set-TextValue 'D33' '168.66'
Set-TextValue 'E33' '  +1.42%  '
Set-TextValue 'D34' '2.92'
Set-TextValue 'E34' '  +6.46%  '
Set-TextValue 'E35' '  -0.22%  '
Set-TextValue 'E36' '  +4.76%  '
Set-TextValue 'D37' '4.81'
Set-TextValue 'E37' '  +2.25%  '
Set-TextValue 'D38' '2.98'
Set-TextValue 'E38' '  +5.64%  '
Set-TextValue 'D39' '3.92'
Set-TextValue 'E39' '  +6.46%  '
Set-TextValue 'E40' '  -0.76%  '
Set-TextValue 'E41' '  +7.54%  '
Set-TextValue 'D42' '99.33'
Set-TextValue 'E42' '  -5.01%  '
Set-TextValue 'E43' '  -0.57%  '
Set-TextValue 'D44' '69.49'
Set-TextValue 'E44' '  -1.82%  '
Set-TextValue 'D45' '12.84'
Set-TextValue 'E45' '  -0.98%  '
Set-TextValue 'E46' '  +0.14%  '
Set-TextValue 'D47' '81.87'
Set-TextValue 'E47' '  +6.22%  '
Set-TextValue 'D48' '111.67'
Set-TextValue 'E48' '  -0.88%  '
Set-TextValue 'E49' '  +4.80%  '
Set-TextValue 'D50' '1.680.88'
Set-TextValue 'E50' '  +1.43%  '
Set-TextValue 'D51' '9.15'
Set-TextValue 'E51' '  +4.09%  '
